$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.669.32"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").Value = "3.799.19"
$ws.Range("E3").Value = "  +0.97%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'596.98"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.62%  "
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +0.64%  "
$ws.Range("E9").Value = "  +1.18%  "
$ws.Range("D10").Value = "'6.30"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.02%  "
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  -0.17%  "
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "4.439.86"
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("D15").Value = "3.777.92"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'18.63"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.46%  "
$ws.Range("D17").Value = "67.696.51"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("E18").Value = "  +2.36%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").Value = "'461.06"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "'9.93"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.50%  "
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("E23").Value = "  -0.08%  "
$ws.Range("D24").Value = "'83.38"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'12.11"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.80%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "3.940.62"
$ws.Range("E29").Value = "  +0.87%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'2.24"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.11%  "
$ws.Range("D32").Value = "'7.29"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").Value = "'29.59"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "
$ws.Range("E34").Value = "  +0.02%  "
$ws.Range("E35").Value = "  -1.10%  "
$ws.Range("E36").Value = "  +0.81%  "
$ws.Range("D37").Value = "'3.38"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.48%  "
$ws.Range("E38").Value = "  +0.45%  "
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("E40").Value = "  +1.06%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D43").Value = "'48.16"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.13%  "
$ws.Range("E44").Value = "  +1.41%  "
$ws.Range("D45").Value = "'43.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.84%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("B47").Value = "ONDO"
$ws.Range("C47").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D47").Value = "'1.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +12.32%  "
$ws.Range("E48").Value = "  +8.56%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'147.40"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.10%  "
$ws.Range("D50").Value = "'394.24"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("E51").Value = "  +0.93%  "
